# Daily attendance processing - 2026-01-13 09:43:14
# For every row in the "Recorded By" column (G), if the value begins with
# "System, " move the leading "System" token to the end of the
# comma-separated list (e.g. "System, a, b" -> "a, b, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val.ToString().StartsWith("System, ")) {
        $rest = $val.ToString().Substring(8)
        $cell.Value = $rest + ", System"
    }
}
